$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells stay text (avoid Excel auto-numeric conversion that would
# strip formatting like trailing zeros or thousand-separator dots).
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "39.493.99"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "2.160.35"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "228.01"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  +2.77%  "
$ws.Range("D7").Value = "63.38"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").Value = "0.0851"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "15.98"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "2.485.08"
$ws.Range("E13").Value = "  +3.04%  "
$ws.Range("D14").Value = "21.98"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "0.808"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "5.49"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "2.157.01"
$ws.Range("D18").Value = "39.545.35"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").Value = "6.17"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "72.08"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "0.0₃0845"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "229.52"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "2.32"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").Value = "9.65"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").Value = "171.93"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "19.83"
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("D30").Value = "1.42"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "2.65"
$ws.Range("E31").Value = "  +5.17%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").Value = "4.58"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D34").Value = "7.02"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "4.68"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").Value = "0.0617"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "2.42"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").Value = "3.63"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "102.28"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").Value = "4.62"
$ws.Range("E41").Value = "  +12.00%  "
$ws.Range("D42").Value = "0.0227"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "17.76"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").Value = "1.523.21"
$ws.Range("D45").Value = "1.20"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.0923"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "7.71"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.368.51"
$ws.Range("E51").Value = "  +3.10%  "

Write-Output "Applied cryptos list update."
